$p = $ppt.ActivePresentation

# Duplicate the previous "Daily Announcements" slide (Lecture 15) to create
# Lecture 16, keeping its layout, placeholders, tags and formatting intact.
$src = $p.Slides.Item(15)
$range = $src.Duplicate()
$ns = $range.Item(1)

# --- Title -----------------------------------------------------------
$ns.Shapes.Item(1).TextFrame.TextRange.Text = "Lecture 16: TUE Oct. 24"

# --- Body --------------------------------------------------------------
$body = $ns.Shapes.Item(2).TextFrame.TextRange
$body.Text = "DMT2 is going to be great! Let" + [char]8217 + "s learn a lot of stuff together!`rProgramming assignment is done (officially).`rWas extended for everyone due to an increase in extension requests.`rHow is it going / did it go?`rQuiz grades have been released!`rSee email for quiz averages, etc.`rWe will go over it today.`rMod. 3 homework is due this Thursday. Don" + [char]8217 + "t put it off!`rToday we begin module 4. "

# Indent the sub-bullets (paragraphs 3, 4, 6, 7) one level deeper.
$body.Paragraphs(3,1).IndentLevel = 2
$body.Paragraphs(4,1).IndentLevel = 2
$body.Paragraphs(6,1).IndentLevel = 2
$body.Paragraphs(7,1).IndentLevel = 2

# Re-apply the bold/italic/underline emphasis on "going to be great".
$emph = $body.Find("going to be great")
$emph.Font.Bold = $true
$emph.Font.Italic = $true
$emph.Font.Underline = $true
